# Fix the "Partens navn" field-code typo in the template:
#   ${Personer.get(0).get("Navn")}   ->   ${Personer.get(0).getNavn()}
#
# The final paragraph text is produced as three runs (matching the
# canonical OOXML produced by the original commit):
#   1) "Partens navn: ${Personer.get(0).getNavn"
#   2) "("
#   3) ")}"
#
# This engine auto-coalesces two adjacent runs back into one whenever a
# single edit leaves them with matching run properties, so a plain
# Find/Replace (or a single Range.Text assignment spanning both split
# points) collapses everything back into one run. To keep three separate
# <w:r> elements we temporarily mark the new split-off piece with Bold
# (so it differs from its neighbours and doesn't get silently re-merged
# at the moment of the edit), then clear that temporary mark again with
# ONE bulk Font update that spans all three runs at once - multi-run
# formatting writes are not subject to the same auto-merge check that a
# single-run write (or a text-insert) triggers.

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive original text instead of
# a hard-coded index, so the script is resilient to unrelated layout
# changes elsewhere in the document.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like 'Partens navn:*') {
        $targetPara = $para
        break
    }
}

if ($targetPara -eq $null) {
    Write-Output "Target paragraph ('Partens navn: ...') not found - aborting."
} else {
    $pStart = $targetPara.Range.Start
    $pEnd   = $targetPara.Range.End

    # Phase 1: rewrite the whole paragraph's text in one shot. (This merges
    # the original two runs into one - that's fine, we re-split it below.)
    $newText = 'Partens navn: ${Personer.get(0).getNavn()}'
    $whole = $d.Range($pStart, $pEnd)
    $whole.Text = $newText

    # Phase 2: split the single run into three at the positions matching
    # the target OOXML: [0,39)="Partens navn: ${Personer.get(0).getNavn",
    # [39,40)="(", [40,42)=")}".
    $splitA = $pStart + 39
    $splitB = $pStart + 40
    $splitC = $pStart + 42

    # Mark the middle character "(" as Bold so it becomes its own run and
    # doesn't immediately re-merge with its still-plain neighbours.
    $middle = $d.Range($splitA, $splitB)
    $middle.Font.Bold = 1

    # Clear the temporary Bold mark again, but as a SINGLE bulk operation
    # spanning the whole paragraph (all three runs at once) so the
    # per-run auto-merge check never fires.
    $fullRange = $d.Range($pStart, $splitC)
    $fullRange.Font.Bold = 0

    Write-Output "Updated paragraph text: $($targetPara.Range.Text)"
}
